# Version 2.4 # 2024-24-11
#
# Table1 worksheet edit:
#   1) Swap the order of the "Present" / "Not present" rows inside the
#      "Comorbidity (%)" block (row 24-25) and fix the capitalisation of
#      "Not present" -> "Not Present".
#   2) Append a new "Falls (%)" block (rows 39-40), following the same
#      two-row layout used by the other binary (%) variables in the table
#      (label + first level on the first row, blank label + second level
#      on the following row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table1")

# --- 1) Comorbidity (%) block: reorder + recase -----------------------
$ws.Cells.Item(24, 1).Value = "Comorbidity (%)"
$ws.Cells.Item(24, 2).Value = "Not Present"
$ws.Cells.Item(24, 3).Value = "   93 (51.7) "

$ws.Cells.Item(25, 1).Value = ""
$ws.Cells.Item(25, 2).Value = "Present"
$ws.Cells.Item(25, 3).Value = "   87 (48.3) "

# --- 2) New Falls (%) block appended at the end of the table -----------
$ws.Cells.Item(39, 1).Value = "Falls (%)"
$ws.Cells.Item(39, 2).Value = "No"
$ws.Cells.Item(39, 3).Value = "  126 (70.0) "

$ws.Cells.Item(40, 1).Value = ""
$ws.Cells.Item(40, 2).Value = "Yes"
$ws.Cells.Item(40, 3).Value = "   54 (30.0) "
